$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.606.28"
$ws.Range("E2").Value = "  -0.64%  "

$ws.Range("D3").Value = "2.667.72"
$ws.Range("E3").Value = "  -0.87%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.614"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.24%  "

$ws.Range("E9").Value = "  +2.19%  "

$ws.Range("E10").Value = "  -0.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.38%  "

$ws.Range("E14").Value = "  -4.39%  "

$ws.Range("D15").Value = "3.147.03"
$ws.Range("E15").Value = "  -0.80%  "

$ws.Range("D16").Value = "65.515.77"
$ws.Range("E16").Value = "  -0.45%  "

$ws.Range("D17").Value = "2.656.33"
$ws.Range("E17").Value = "  -1.18%  "

$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("E19").Value = "  -1.73%  "

$ws.Range("E20").Value = "  -2.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.66%  "

$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.62%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.73%  "

$ws.Range("E26").Value = "  +2.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.54%  "

$ws.Range("E28").Value = "  -3.76%  "

$ws.Range("E29").Value = "  -3.17%  "

$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("E31").Value = "  -3.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "530.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.13%  "

$ws.Range("E33").Value = "  +0.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.49%  "

$ws.Range("E35").Value = "  -0.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.422"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "156.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.99%  "

$ws.Range("E40").Value = "  -3.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "162.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.00%  "

$ws.Range("E43").Value = "  -1.32%  "

$ws.Range("E44").Value = "  +1.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0608"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.55%  "

$ws.Range("E47").Value = "  -2.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0257"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.22%  "

$ws.Range("D49").Value = "0.0₆0253"
$ws.Range("E49").Value = "  +6.59%  "

$ws.Range("E50").Value = "  -1.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.36%  "
